$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark: it used to mark the spot right after the
#    page-break run (just before "DER:"); the latest edit happened instead
#    mid-word inside the "Tomando en cuenta..." paragraph, so remove it from
#    the old spot and re-add it at the new one.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Extend the "Se decidió separar cada funcionalidad..." bullet with the
#    new trailing sentence about being more specific.
# ---------------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("Se decidió separar cada funcionalidad en ítems para tener el mayor control de quién tiene acceso a qué pantalla.")
$target.Collapse(0)
$target.InsertAfter(" Teniendo en cuenta que en un futuro se puede querer descentralizar las tareas del administrador y se desconoce las funcionalidades a realizar por cada rol, se prefirió ser más específico.")

# ---------------------------------------------------------------------------
# 3. Re-insert "_GoBack" at its new location: right after "...una operac",
#    in the middle of the word "operación" inside the next bullet. No text
#    actually changes there -- only the bookmark moves.
# ---------------------------------------------------------------------------
$goback = $d.Content
$goback.Find.Execute("Tomando en cuenta que no tenemos una operac")
$goback.Collapse(0)
$d.Bookmarks.Add("_GoBack", $goback)
